# Auto-generated edit script: adds 2023-11-18 crime data
# Updates per-sheet 2023 (column J) figures and a couple of 2022 (column I) corrections
$wb = $excel.ActiveWorkbook

$sheetEdits = @{
    'Citywide Totals' = @{ "J2" = 6801; "I3" = 7491; "J3" = 7195; "J4" = 1569; "J5" = 565; "J6" = 9590; "I7" = 26230; "J7" = 25720 }
    'Logan Square' = @{ "J2" = 66; "J6" = 245; "J7" = 371 }
    'Austin' = @{ "J4" = 87; "J6" = 571; "J7" = 1615 }
    'South Chicago' = @{ "J3" = 191; "J7" = 513 }
    'Garfield Park' = @{ "J3" = 387; "J5" = 49; "J6" = 404; "J7" = 1162 }
    'West Pullman' = @{ "J2" = 128; "J3" = 135; "J4" = 18; "J7" = 375 }
    'Grand Crossing' = @{ "J2" = 235; "J3" = 267; "J6" = 232; "J7" = 795 }
    'New City' = @{ "J3" = 184; "J7" = 644 }
    'Woodlawn' = @{ "J3" = 159; "J7" = 395 }
    'By Neighborhood' = @{ "J2" = 203; "J7" = 745; "J8" = 1615; "J10" = 188; "J11" = 443; "J15" = 305; "J19" = 751; "J20" = 540; "J25" = 129; "J28" = 8; "I29" = 1556; "J29" = 1393; "J31" = 258; "J32" = 42; "J33" = 1162; "J37" = 795; "J42" = 1105; "J43" = 220; "J47" = 191; "J48" = 293; "J50" = 153; "J51" = 313; "J52" = 647; "J53" = 371; "J54" = 494; "J55" = 394; "J63" = 81; "J65" = 644; "J67" = 966; "J73" = 248; "J74" = 29; "J79" = 724; "J83" = 513; "J85" = 1063; "J88" = 270; "J91" = 297; "J94" = 276; "J95" = 375; "J97" = 236; "J99" = 395; "I101" = 26230; "J101" = 25720 }
    'Gage Park' = @{ "J2" = 91; "J6" = 84; "J7" = 258 }
    'North Lawndale' = @{ "J2" = 246; "J3" = 357; "J7" = 966 }
    'Loop' = @{ "J2" = 121; "J6" = 231; "J7" = 494 }
    'Englewood' = @{ "I3" = 529; "J3" = 495; "J6" = 353; "I7" = 1556; "J7" = 1393 }
    'Lake View' = @{ "J3" = 54; "J4" = 46; "J7" = 293 }
    'Chatham' = @{ "J6" = 291; "J7" = 751 }
    'Humboldt Park' = @{ "J3" = 220; "J7" = 1105 }
    'Avondale' = @{ "J6" = 107; "J7" = 188 }
    'Lower West Side' = @{ "J2" = 77; "J7" = 394 }
    'Washington Park' = @{ "J6" = 75; "J7" = 297 }
    'Roseland' = @{ "J2" = 201; "J3" = 244; "J5" = 20; "J7" = 724 }
    'Chicago Lawn' = @{ "J3" = 184; "J6" = 153; "J7" = 540 }
    'Auburn Gresham' = @{ "J2" = 230; "J6" = 239; "J7" = 745 }
    'West Loop' = @{ "J4" = 22; "J7" = 276 }
    'East Side' = @{ "J6" = 23; "J7" = 129 }
    'Kenwood' = @{ "J6" = 89; "J7" = 191 }
    'Brighton Park' = @{ "J2" = 86; "J7" = 305 }
    'Lincoln Square' = @{ "J6" = 53; "J7" = 153 }
    'Belmont Cragin' = @{ "J2" = 129; "J3" = 80; "J7" = 443 }
    'Portage Park' = @{ "J6" = 88; "J7" = 248 }
    'Albany Park' = @{ "J2" = 60; "J7" = 203 }
    'West Town' = @{ "J5" = 2; "J7" = 236 }
    'United Center' = @{ "J6" = 137; "J7" = 270 }
    'Galewood' = @{ "J6" = 19; "J7" = 42 }
    'Little Italy, UIC' = @{ "J6" = 125; "J7" = 313 }
    'Hyde Park' = @{ "J6" = 131; "J7" = 220 }
    'South Shore' = @{ "J2" = 287; "J3" = 372; "J6" = 309; "J7" = 1063 }
    'Little Village' = @{ "J3" = 184; "J7" = 647 }
    'Printers Row' = @{ "J2" = 5; "J6" = 13; "J7" = 29 }
    'Edison Park' = @{ "J6" = 4; "J7" = 8 }
}

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $sheetEdits[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
